# Update countries & provincias Spain
# Applies updated COVID case stats for "Estados Unidos" (row 4) and
# "Alemania" (row 9), plus refreshed stats for Cabo Verde which now
# overtakes Liberia and Birmania in the ranked list (rows 138-140).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 1185167   # B4 Casos totales
$ws.Cells.Item(4, 3).Value = 24393     # C4 Nuevos casos
$ws.Cells.Item(4, 5).Value = 938453    # E4 Recuperados
$ws.Cells.Item(4, 7).Value = 1051      # G4 Casos criticos
$ws.Cells.Item(4, 8).Value = 68495     # H4 Muertes

# --- Row 9: Alemania ---
$ws.Cells.Item(9, 2).Value = 165565    # B9 Casos totales
$ws.Cells.Item(9, 3).Value = 598       # C9 Nuevos casos
$ws.Cells.Item(9, 5).Value = 28117     # E9 Recuperados
$ws.Cells.Item(9, 7).Value = 36        # G9 Casos criticos
$ws.Cells.Item(9, 8).Value = 6848      # H9 Muertes

# --- Rows 138-140: Cabo Verde overtakes Liberia and Birmania ---
# Row 138 becomes Cabo Verde with refreshed data
$ws.Cells.Item(138, 1).Value = "Cabo Verde"
$ws.Cells.Item(138, 2).Value = 165
$ws.Cells.Item(138, 3).Value = 13
$ws.Cells.Item(138, 4).Value = 33
$ws.Cells.Item(138, 5).Value = 130
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 2

# Row 139 becomes Liberia (previously row 138's data)
$ws.Cells.Item(139, 1).Value = "Liberia"
$ws.Cells.Item(139, 2).Value = 158
$ws.Cells.Item(139, 3).Value = 4
$ws.Cells.Item(139, 4).Value = 58
$ws.Cells.Item(139, 5).Value = 82
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 18

# Row 140 becomes Birmania (previously row 139's data)
$ws.Cells.Item(140, 1).Value = "Birmania"
$ws.Cells.Item(140, 2).Value = 155
$ws.Cells.Item(140, 3).Value = 4
$ws.Cells.Item(140, 4).Value = 43
$ws.Cells.Item(140, 5).Value = 106
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 6
